$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# The slide contains two identical confusion-matrix tables (shapes 1 and 2,
# labelled "WIKI" and "SKLEARN" respectively by nearby textboxes). The
# second table's header cell currently reads "Ground truth" and must be
# changed to "Predicted" (sklearn's confusion_matrix columns are the
# predicted labels, not ground truth -- that's the "blind spot").

$table2 = $s.Shapes.Item(2).Table
$table2.Cell(1, 2).Shape.TextFrame.TextRange.Text = "Predicted"
